$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 26.666666
$ws.Range("I8").Value = 14.25
$ws.Range("J8").Value = 51.5
$ws.Range("K8").Value = 42.75
$ws.Range("L8").Value = 154.5
$ws.Range("M8").Value = 96.25
$ws.Range("N8").Value = -432.5
$ws.Range("H80").Value = 2007.375
$ws.Range("I80").Value = 1446
$ws.Range("J80").Value = 2568.75
$ws.Range("K80").Value = 4338
$ws.Range("L80").Value = 7706.25
$ws.Range("M80").Value = -3340
$ws.Range("N80").Value = -9702.25
$ws.Range("H83").Value = 2007.375
$ws.Range("I83").Value = 1446
$ws.Range("J83").Value = 2568.75
$ws.Range("K83").Value = 13014
$ws.Range("L83").Value = 23118.75
$ws.Range("M83").Value = -8022
$ws.Range("N83").Value = -33102.75
$ws.Range("H132").Value = 6000
$ws.Range("I132").Value = 6000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 18000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -15470
$ws.Range("H138").Value = 3361.5
$ws.Range("I138").Value = 3937
$ws.Range("J138").Value = 3233.611
$ws.Range("K138").Value = 11811
$ws.Range("L138").Value = 9700.832999999999
$ws.Range("M138").Value = -6671
$ws.Range("N138").Value = -19980.833

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 29929.334
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 29929.334
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 29929.334
$ws.Range("N76").Value = -30605.334
$ws.Range("H79").Value = 29929.334
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 29929.334
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 29929.334
$ws.Range("N79").Value = -32269.334
$ws.Range("H81").Value = 30000
$ws.Range("I81").Value = 30000
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 30000
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -29002
$ws.Range("H84").Value = 30000
$ws.Range("I84").Value = 30000
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 90000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -85008
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H20").Value = 1749.5
$ws.Range("I20").Value = 500
$ws.Range("J20").Value = 2999
$ws.Range("K20").Value = 500
$ws.Range("L20").Value = 2999
$ws.Range("M20").Value = -253
$ws.Range("N20").Value = -3493

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 191.33333
$ws.Range("I22").Value = 191.33333
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 191.33333
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -18.33332999999999
$ws.Range("H86").Value = 4201
$ws.Range("I86").Value = 5163.75
$ws.Range("J86").Value = 350
$ws.Range("K86").Value = 5163.75
$ws.Range("L86").Value = 350
$ws.Range("M86").Value = -4040.75
$ws.Range("N86").Value = -2596
$ws.Range("H89").Value = 4201
$ws.Range("I89").Value = 5163.75
$ws.Range("J89").Value = 350
$ws.Range("K89").Value = 25818.75
$ws.Range("L89").Value = 1750
$ws.Range("M89").Value = -20202.75
$ws.Range("N89").Value = -12982
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H99").Value = 1011
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 1011
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 1011
$ws.Range("N99").Value = -4007
$ws.Range("H100").Value = 13374
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 13374
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 13374
$ws.Range("N100").Value = -15538
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H134").Value = 2391.1667
$ws.Range("I134").Value = 1536.2
$ws.Range("J134").Value = 6666
$ws.Range("K134").Value = 4608.6
$ws.Range("L134").Value = 19998
$ws.Range("M134").Value = -2073.6
$ws.Range("N134").Value = -25068

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 50.6
$ws.Range("I7").Value = 52.333332
$ws.Range("J7").Value = 35
$ws.Range("K7").Value = 52.333332
$ws.Range("L7").Value = 35
$ws.Range("M7").Value = 60.666668
$ws.Range("N7").Value = -261
$ws.Range("H16").Value = 3799.5
$ws.Range("I16").Value = 3799
$ws.Range("J16").Value = 3800
$ws.Range("K16").Value = 3799
$ws.Range("L16").Value = 3800
$ws.Range("M16").Value = -3512
$ws.Range("N16").Value = -4374
$ws.Range("H22").Value = 643
$ws.Range("I22").Value = 96.333336
$ws.Range("J22").Value = 1189.6666
$ws.Range("K22").Value = 96.333336
$ws.Range("L22").Value = 1189.6666
$ws.Range("M22").Value = 253.666664
$ws.Range("N22").Value = -1889.6666
$ws.Range("H47").Value = 19000
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 19000
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 19000
$ws.Range("N47").Value = -20132
$ws.Range("H68").Value = 47031.668
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 47031.668
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 47031.668
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -48529.668
$ws.Range("H71").Value = 47031.668
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 47031.668
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 141095.004
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -148583.004
$ws.Range("H95").Value = 10623.5
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 10623.5
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 10623.5
$ws.Range("N95").Value = -16115.5
$ws.Range("H105").Value = 2999.25
$ws.Range("I105").Value = 2999.25
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2999.25
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -1252.25
$ws.Range("H107").Value = 1371.8
$ws.Range("I107").Value = 286.66666
$ws.Range("J107").Value = 2999.5
$ws.Range("K107").Value = 286.66666
$ws.Range("L107").Value = 2999.5
$ws.Range("M107").Value = 1633.33334
$ws.Range("N107").Value = -6839.5
$ws.Range("H113").Value = 3799.5
$ws.Range("I113").Value = 3799
$ws.Range("J113").Value = 3800
$ws.Range("K113").Value = 3799
$ws.Range("L113").Value = 3800
$ws.Range("M113").Value = -1629
$ws.Range("N113").Value = -8140
$ws.Range("H122").Value = 3800
$ws.Range("I122").Value = 3800
$ws.Range("J122").Value = 3800
$ws.Range("K122").Value = 11400
$ws.Range("L122").Value = 11400
$ws.Range("M122").Value = -8950
$ws.Range("N122").Value = -16300

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H103").Value = 3399.5
$ws.Range("I103").Value = 799.6667
$ws.Range("J103").Value = 4959.4
$ws.Range("K103").Value = 2399.0001
$ws.Range("L103").Value = 14878.2
$ws.Range("M103").Value = -1520.0001
$ws.Range("N103").Value = -16636.2
$ws.Range("H140").Value = 1442.6666
$ws.Range("I140").Value = 1442.6666
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 4327.9998
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 852.0002000000004

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 3215202.5
$ws.Range("I3").Value = 5000037.5
$ws.Range("J3").Value = 2501268.5
$ws.Range("K3").Value = 5000037.5
$ws.Range("L3").Value = 2501268.5
$ws.Range("M3").Value = -4999921.5
$ws.Range("N3").Value = -2501500.5
$ws.Range("H97").Value = 371.75
$ws.Range("I97").Value = 371.75
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 371.75
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 124.25
$ws.Range("N97").ClearContents()
$ws.Range("A98").Value = None
$ws.Range("H98").Value = 9302
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 9302
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 9302
$ws.Range("N98").Value = -15292
$ws.Range("H105").Value = 18434.2
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 18434.2
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 18434.2
$ws.Range("N105").Value = -25422.2
$ws.Range("H113").Value = 2636.5
$ws.Range("I113").Value = 2498.5
$ws.Range("J113").Value = 2774.5
$ws.Range("K113").Value = 2498.5
$ws.Range("L113").Value = 2774.5
$ws.Range("M113").Value = -328.5
$ws.Range("N113").Value = -7114.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1324.8334
$ws.Range("I22").Value = 581
$ws.Range("J22").Value = 2812.5
$ws.Range("K22").Value = 581
$ws.Range("L22").Value = 2812.5
$ws.Range("M22").Value = -286
$ws.Range("N22").Value = -3402.5
$ws.Range("H27").Value = 1324.8334
$ws.Range("I27").Value = 581
$ws.Range("J27").Value = 2812.5
$ws.Range("K27").Value = 581
$ws.Range("L27").Value = 2812.5
$ws.Range("M27").Value = -474
$ws.Range("N27").Value = -3026.5
$ws.Range("H68").Value = 1715.8
$ws.Range("I68").Value = 1715.8
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1715.8
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -966.8
$ws.Range("H71").Value = 1715.8
$ws.Range("I71").Value = 1715.8
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 8579
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -4835

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 19666.334
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 19666.334
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 19666.334
$ws.Range("N69").Value = -21164.334
$ws.Range("H72").Value = 19666.334
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 19666.334
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 58999.00199999999
$ws.Range("N72").Value = -66487.00199999999
$ws.Range("H82").Value = 39995
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 39995
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 39995
$ws.Range("N82").Value = -40761
$ws.Range("H85").Value = 39995
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 39995
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 39995
$ws.Range("N85").Value = -42647
$ws.Range("H92").Value = 24000
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 24000
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 24000
$ws.Range("N92").Value = -28992
$ws.Range("H95").Value = 27500
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 27500
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 27500
$ws.Range("N95").Value = -32992
$ws.Range("H96").Value = 2151
$ws.Range("I96").Value = 2151
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 2151
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -778
$ws.Range("H97").Value = 11838.667
$ws.Range("I97").Value = 10520
$ws.Range("J97").Value = 12498
$ws.Range("K97").Value = 10520
$ws.Range("L97").Value = 12498
$ws.Range("M97").Value = -9529
$ws.Range("N97").Value = -14480
$ws.Range("H107").Value = 1252.7858
$ws.Range("I107").Value = 1440.1666
$ws.Range("J107").Value = 1112.25
$ws.Range("K107").Value = 4320.4998
$ws.Range("L107").Value = 3336.75
$ws.Range("M107").Value = -2400.4998
$ws.Range("N107").Value = -7176.75
